# Chap05: Only 2 pictures remaining: pumping exp schema and probe away
# from the dots. Caption to do in today's pictures.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The two formula shapes ("Rectangle 25", ids 35 and 39) live inside the
# nested group "Groupe 3" (top-level shape 2). Locate shapes by Id/Name
# instead of a hard-coded collection index so the script is resilient to
# ordering.
$g = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Groupe 3") {
        $g = $cand
        break
    }
}

$sh35 = $null
$sh39 = $null
for ($i = 1; $i -le $g.GroupItems.Count; $i++) {
    $cand = $g.GroupItems.Item($i)
    if ($cand.Id -eq 35) { $sh35 = $cand }
    if ($cand.Id -eq 39) { $sh39 = $cand }
}

# -----------------------------------------------------------------
# Shape "Rectangle 25" (id 35) at off (5731927,4850197): Sz=-1 -> Sz= +1
# -----------------------------------------------------------------
$tr35 = $sh35.TextFrame.TextRange

# Edit right-to-left so earlier inserts don't shift later offsets.
$tr35.Characters(4,1).Text = "+"
$tr35.Characters(3,1).Text = "= "

# Resize the shape (cx 531848->562244, cy 256217->366364 EMU). The
# Width/Height setters store the raw point value * 12700 as the EMU
# extent, so convert directly (a tiny epsilon corrects float truncation).
$sh35.Width = 44.27123110236221
$sh35.Height = 28.847609055118113

# -----------------------------------------------------------------
# Shape "Rectangle 25" (id 39) at off (4677026,4850196): Sz=+1 -> Sz= -1 (en dash)
# -----------------------------------------------------------------
$tr39 = $sh39.TextFrame.TextRange

$tr39.Characters(4,1).Text = "–"
$tr39.Characters(3,1).Text = "= "

# Resize the shape (cx 637595->535832, cy 301807->366364 EMU).
$sh39.Width = 42.191546062992124
$sh39.Height = 28.847609055118113

# -----------------------------------------------------------------
# Remove the stray top-level "Freeform 88" shape (id 54). A second,
# unrelated "Freeform 88" (id 58) lives inside the group and must stay.
# -----------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Freeform 88") {
        $sh.Delete()
        break
    }
}
